# Update for March 16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ontario")

# Correction to March 15 case count
$ws.Range("B21").Value = 146

# New data row for March 16, 2020
$ws.Range("A22").Value = "2020-03-16"
$ws.Range("A22").NumberFormat = $ws.Range("A21").NumberFormat
$ws.Range("B22").Value = 177
$ws.Range("C22").Value = 10178

$ws.Range("D22").Formula = "=A22-`$A`$2"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Formula = "=D22-D21"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("F22").Formula = "=(B22/B21)^(1/E22)-1"
$ws.Range("F22").NumberFormat = $ws.Range("F21").NumberFormat

$ws.Range("G22").Value = "https://www.cbc.ca/news/canada/toronto/ontario-covid-19-coronavirus-monday-1.5498849"
$ws.Hyperlinks.Add($ws.Range("G22"), "https://www.cbc.ca/news/canada/toronto/ontario-covid-19-coronavirus-monday-1.5498849") | Out-Null
$ws.Range("G22").Style = $ws.Range("G21").Style

$wb.Save()
